$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.495.90"
$ws.Range("E2").Value = "'  -3.28%  "

# Row 3
$ws.Range("D3").Value = "'1.996.79"
$ws.Range("E3").Value = "'  -5.94%  "

# Row 4
$ws.Range("E4").Value = "'  +0.07%  "

# Row 5
$ws.Range("D5").Value = "'329.66"
$ws.Range("E5").Value = "'  -4.89%  "

# Row 7
$ws.Range("D7").Value = "'0.5010"
$ws.Range("E7").Value = "'  -4.30%  "

# Row 8
$ws.Range("D8").Value = "'0.4238"
$ws.Range("E8").Value = "'  -5.03%  "

# Row 9
$ws.Range("D9").Value = "'53.28"
$ws.Range("E9").Value = "'  -2.09%  "

# Row 10
$ws.Range("D10").Value = "'0.08933"
$ws.Range("E10").Value = "'  -4.90%  "

# Row 11
$ws.Range("E11").Value = "'  -4.69%  "

# Row 12
$ws.Range("D12").Value = "'23.28"
$ws.Range("E12").Value = "'  -7.67%  "

# Row 13
$ws.Range("D13").Value = "'8.101"
$ws.Range("E13").Value = "'  -6.87%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "'6.522"
$ws.Range("E14").Value = "'  -6.28%  "

# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "'1.931.87"
$ws.Range("E15").Value = "'  -9.13%  "

# Row 16
$ws.Range("D16").Value = "'95.72"
$ws.Range("E16").Value = "'  -6.39%  "

# Row 18
$ws.Range("D18").Value = "'0.00001111"
$ws.Range("E18").Value = "'  -4.76%  "

# Row 19
$ws.Range("D19").Value = "'0.06630"
$ws.Range("E19").Value = "'  -1.46%  "

# Row 20
$ws.Range("D20").Value = "'19.71"
$ws.Range("E20").Value = "'  -7.97%  "

# Row 21
$ws.Range("D21").Value = "'1.008"
$ws.Range("E21").Value = "'  +0.12%  "

# Row 22
$ws.Range("D22").Value = "'5.972"
$ws.Range("E22").Value = "'  -5.81%  "

# Row 23
$ws.Range("D23").Value = "'29.505.12"
$ws.Range("E23").Value = "'  -3.37%  "

# Row 24
$ws.Range("D24").Value = "'11.94"
$ws.Range("E24").Value = "'  -6.30%  "

# Row 25
$ws.Range("D25").Value = "'2.255"
$ws.Range("E25").Value = "'  -3.18%  "

# Row 26
$ws.Range("E26").Value = "'  -2.55%  "

# Row 27
$ws.Range("D27").Value = "'20.67"
$ws.Range("E27").Value = "'  -6.87%  "

# Row 28
$ws.Range("D28").Value = "'6.568"
$ws.Range("E28").Value = "'  -4.51%  "

# Row 29
$ws.Range("D29").Value = "'2.334"
$ws.Range("E29").Value = "'  -8.34%  "

# Row 30
$ws.Range("D30").Value = "'127.79"
$ws.Range("E30").Value = "'  -4.77%  "

# Row 31
$ws.Range("D31").Value = "'1.049"
$ws.Range("E31").Value = "'  -9.47%  "

# Row 32
$ws.Range("D32").Value = "'0.09950"
$ws.Range("E32").Value = "'  -6.25%  "

# Row 33
$ws.Range("D33").Value = "'1.574"
$ws.Range("E33").Value = "'  -11.36%  "

# Row 34
$ws.Range("D34").Value = "'5.855"
$ws.Range("E34").Value = "'  -6.84%  "

# Row 35
$ws.Range("E35").Value = "'  -4.52%  "

# Row 36
$ws.Range("D36").Value = "'9.581"
$ws.Range("E36").Value = "'  -9.76%  "

# Row 37
$ws.Range("D37").Value = "'0.02467"
$ws.Range("E37").Value = "'  -7.39%  "

# Row 38
$ws.Range("D38").Value = "'0.06343"
$ws.Range("E38").Value = "'  -7.55%  "

# Row 39
$ws.Range("D39").Value = "'1.289"
$ws.Range("E39").Value = "'  -3.86%  "

# Row 40
$ws.Range("D40").Value = "'0.6527"
$ws.Range("E40").Value = "'  -8.43%  "

# Row 41
$ws.Range("D41").Value = "'11.71"
$ws.Range("E41").Value = "'  -7.68%  "

# Row 42
$ws.Range("D42").Value = "'0.2067"
$ws.Range("E42").Value = "'  -7.78%  "

# Row 44
$ws.Range("D44").Value = "'0.6334"
$ws.Range("E44").Value = "'  -8.70%  "

# Row 45
$ws.Range("D45").Value = "'2.209"
$ws.Range("E45").Value = "'  -7.28%  "

# Row 46
$ws.Range("D46").Value = "'13.41"
$ws.Range("E46").Value = "'  -8.66%  "

# Row 47
$ws.Range("D47").Value = "'1.273"
$ws.Range("E47").Value = "'  -4.37%  "

# Row 48
$ws.Range("E48").Value = "'  -3.41%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.06997"
$ws.Range("E49").Value = "'  -3.35%  "

# Row 50
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "'0.00000000324"
$ws.Range("E50").Value = "'  -5.67%  "

# Row 51
$ws.Range("E51").Value = "'  -5.05%  "
